$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.375.88'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.712.30'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.91%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +2.77%  '
$ws.Range('D9').Value = '2.711.96'
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.146'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '3.207.24'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '68.439.55'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('D18').Value = '2.722.83'
$ws.Range('E18').Value = '  +2.74%  '
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '369.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('E22').Value = '  +1.29%  '
$ws.Range('E23').Value = '  +3.00%  '
$ws.Range('E24').Value = '  -0.87%  '
$ws.Range('E25').Value = '  -2.66%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '579.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('E34').Value = '  +5.81%  '
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '160.45'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.379'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.16%  '
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.87'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').Value = '  -1.93%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '0.0₆0311'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.04%  '
$ws.Range('E48').Value = '  +3.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '155.08'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('E51').Value = '  +3.51%  '
